$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 220
$ws.Range("I33").Value = 174
$ws.Range("J33").Value = 756.6667
$ws.Range("K33").Value = 174
$ws.Range("L33").Value = 756.6667
$ws.Range("M33").Value = 55
$ws.Range("N33").Value = -1214.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 2830.3333
$ws.Range("I100").Value = 1665.2858
$ws.Range("K100").Value = 1665.2858
$ws.Range("M100").Value = -1124.2858

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H115").Value = 1852.4
$ws.Range("I115").Value = 1852.4
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 5557.200000000001
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = -3990.200000000001
$ws.Range("N115").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 4948.9644
$ws.Range("I137").Value = 1154.7142
$ws.Range("K137").Value = 3464.1426
$ws.Range("M137").Value = -914.1425999999997

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1603.94
$ws.Range("I138").Value = 814.625
$ws.Range("J138").Value = 1975.3823
$ws.Range("K138").Value = 2443.875
$ws.Range("L138").Value = 5926.1469
$ws.Range("M138").Value = 2696.125
$ws.Range("N138").Value = -16206.1469

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2129.125
$ws.Range("I141").Value = 1214
$ws.Range("J141").Value = 8535
$ws.Range("K141").Value = 3642
$ws.Range("L141").Value = 25605
$ws.Range("M141").Value = 1538
$ws.Range("N141").Value = -35965

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 325
$ws.Range("I4").Value = 325
$ws.Range("K4").Value = 325
$ws.Range("M4").Value = -209

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 236.84616
$ws.Range("I5").Value = 270.9
$ws.Range("J5").Value = 123.333336
$ws.Range("K5").Value = 270.9
$ws.Range("L5").Value = 123.333336
$ws.Range("M5").Value = -158.9
$ws.Range("N5").Value = -347.333336

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 66336
$ws.Range("J12").Value = 66336
$ws.Range("L12").Value = 66336
$ws.Range("N12").Value = -66682

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 1000000
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7911.8413
$ws.Range("I32").Value = 7195.575
$ws.Range("J32").Value = 13721.556
$ws.Range("K32").Value = 7195.575
$ws.Range("L32").Value = 13721.556
$ws.Range("M32").Value = -6908.575
$ws.Range("N32").Value = -14295.556

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2259.524
$ws.Range("I45").Value = 2057.5625
$ws.Range("K45").Value = 2057.5625
$ws.Range("M45").Value = -1680.5625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 43096.5
$ws.Range("J139").Value = 43096.5
$ws.Range("L139").Value = 43096.5
$ws.Range("N139").Value = -53376.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 236.84616
$ws.Range("I4").Value = 270.9
$ws.Range("J4").Value = 123.333336
$ws.Range("K4").Value = 270.9
$ws.Range("L4").Value = 123.333336
$ws.Range("M4").Value = -155.9
$ws.Range("N4").Value = -353.333336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4364.4155
$ws.Range("I31").Value = 2121.875
$ws.Range("J31").Value = 4855.9316
$ws.Range("K31").Value = 2121.875
$ws.Range("L31").Value = 4855.9316
$ws.Range("M31").Value = -1826.875
$ws.Range("N31").Value = -5445.9316

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4364.4155
$ws.Range("I34").Value = 2121.875
$ws.Range("J34").Value = 4855.9316
$ws.Range("K34").Value = 2121.875
$ws.Range("L34").Value = 4855.9316
$ws.Range("M34").Value = -1919.875
$ws.Range("N34").Value = -5259.9316

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 446780.88
$ws.Range("I140").Value = 573137.3
$ws.Range("J140").Value = 4533.5
$ws.Range("K140").Value = 1719411.9
$ws.Range("L140").Value = 13600.5
$ws.Range("M140").Value = -1714231.9
$ws.Range("N140").Value = -23960.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 100004296
$ws.Range("I141").Value = 166670020
$ws.Range("J141").Value = 5724.75
$ws.Range("K141").Value = 500010060
$ws.Range("L141").Value = 17174.25
$ws.Range("M141").Value = -500004880
$ws.Range("N141").Value = -27534.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 2008105.6
$ws.Range("I2").Value = 6024116
$ws.Range("J2").Value = 100.5
$ws.Range("K2").Value = 6024116
$ws.Range("L2").Value = 100.5
$ws.Range("M2").Value = -6024003
$ws.Range("N2").Value = -326.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1621.5
$ws.Range("I113").Value = 1619.8
$ws.Range("J113").Value = 1625.75
$ws.Range("K113").Value = 1619.8
$ws.Range("L113").Value = 1625.75
$ws.Range("M113").Value = 550.2
$ws.Range("N113").Value = -5965.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 628.5
$ws.Range("I22").Value = 517.46155
$ws.Range("J22").Value = 788.8889
$ws.Range("K22").Value = 517.46155
$ws.Range("L22").Value = 788.8889
$ws.Range("M22").Value = -222.46155
$ws.Range("N22").Value = -1378.8889

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 628.5
$ws.Range("I27").Value = 517.46155
$ws.Range("J27").Value = 788.8889
$ws.Range("K27").Value = 517.46155
$ws.Range("L27").Value = 788.8889
$ws.Range("M27").Value = -410.46155
$ws.Range("N27").Value = -1002.8889

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 127237.75
$ws.Range("I122").Value = 168817.33
$ws.Range("J122").Value = 2499
$ws.Range("K122").Value = 506451.99
$ws.Range("L122").Value = 7497
$ws.Range("M122").Value = -504001.99
$ws.Range("N122").Value = -12397

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1964.9056
$ws.Range("I132").Value = 1642.5209
$ws.Range("J132").Value = 5059.8
$ws.Range("K132").Value = 4927.5627
$ws.Range("L132").Value = 15179.4
$ws.Range("M132").Value = -2397.5627
$ws.Range("N132").Value = -20239.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 720.6667
$ws.Range("I107").Value = 575.3333
$ws.Range("J107").Value = 866
$ws.Range("K107").Value = 1725.9999
$ws.Range("L107").Value = 2598
$ws.Range("M107").Value = 194.0001
$ws.Range("N107").Value = -6438

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H140").Value = 52285.125
$ws.Range("J140").Value = 52285.125
$ws.Range("L140").Value = 52285.125
$ws.Range("N140").Value = -62645.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H141").Value = 38499.832
$ws.Range("J141").Value = 38499.832
$ws.Range("L141").Value = 38499.832
$ws.Range("N141").Value = -48859.832
